$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New single-deck bus rows (42-46)
$newBuses = @(
    @{ Row = 42; Name = "Northern General SE6"; Year = 1933; Order = 1; Speed = 40; Capacity = 44 },
    @{ Row = 43; Name = "Leyland Tiger TS";      Year = 1927; Order = 1; Speed = 42; Capacity = 35 },
    @{ Row = 44; Name = "Leyland National";      Year = 1972; Order = 1; Speed = 52; Capacity = 58 },
    @{ Row = 45; Name = "Dennis Dart";           Year = 1989; Order = 1; Speed = 56; Capacity = 50 },
    @{ Row = 46; Name = "Sentinel Steam Bus";    Year = 1924; Order = 1; Speed = 36; Capacity = 32 }
)

foreach ($bus in $newBuses) {
    $r = $bus.Row
    $ws.Range("A$r").Value = $bus.Name
    $ws.Range("B$r").Value = $bus.Year
    $ws.Range("C$r").Value = $bus.Order
    $ws.Range("D$r").Value = "Bus"
    $ws.Range("E$r").Formula = "=IF(B$r > 1900, ((B$r-1900)*10)+400+C$r, ((B$r-1730)*2)+C$r)+VLOOKUP(D$r,'ID Scheme'!`$A`$2:`$B`$4,2)"
    $ws.Range("F$r").Value = $bus.Speed
    $ws.Range("G$r").Value = $bus.Capacity
    $ws.Range("H$r").Formula = "=SQRT(F$r*G$r)/`$B`$1"
    $ws.Range("I$r").Formula = "=H$r*0.9"
    $ws.Range("J$r").Value = "x"
    $ws.Range("H${r}:J${r}").NumberFormat = "0"
}

# Freeze the header row (row 3) so the pane splits below it, then scroll down
# and leave the new bottom-right entry selected, matching the saved view.
$ws.Activate()
$ws.Range("A4").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("J42").Select()
